$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (shifts existing D:K quarterly data to F:M)
$ws.Range("D:E").Insert()

# The new D:E columns land with default formatting; copy number/date formats
# from the (now shifted) F:G columns, which still carry the correct per-row styles.
# Scoped to the three contiguous data blocks so header-only / blank separator
# rows (5,6,36,37,78,79) are not touched.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

# Write the two new quarters of data into D:E, plus the handful of
# shifted-column corrections called out in the source data.
$data = @(
    @(7, "D", 43465),
    @(7, "E", 43373),
    @(8, "D", 13500),
    @(8, "E", 12500),
    @(9, "D", "NA"),
    @(9, "E", "NA"),
    @(10, "D", "NA"),
    @(10, "E", "NA"),
    @(12, "D", "NA"),
    @(12, "E", "NA"),
    @(13, "D", 0),
    @(13, "E", 0),
    @(14, "D", 0),
    @(14, "E", 0),
    @(15, "D", -200),
    @(15, "E", -100),
    @(17, "D", 2800),
    @(17, "E", 2000),
    @(18, "D", 10700),
    @(18, "E", 10500),
    @(20, "D", -5400),
    @(20, "E", -4900),
    @(21, "D", 5700),
    @(21, "E", 6100),
    @(22, "D", 0),
    @(22, "E", 0),
    @(23, "D", 5300),
    @(23, "E", 5700),
    @(24, "D", 1100),
    @(24, "E", 1200),
    @(25, "D", 0),
    @(25, "E", 0),
    @(26, "D", 4200),
    @(26, "E", 4500),
    @(27, "D", 4200),
    @(27, "E", 4500),
    @(28, "D", 0),
    @(28, "E", 0),
    @(29, "D", 0),
    @(29, "E", "NA"),
    @(30, "D", 0),
    @(30, "E", 0),
    @(31, "D", 0),
    @(31, "E", 0),
    @(32, "D", 5400),
    @(32, "E", 4900),
    @(33, "D", 4200),
    @(33, "E", 4500),
    @(34, "D", 0),
    @(34, "E", 0),
    @(35, "D", 4200),
    @(35, "E", 4500),
    @(38, "D", 43465),
    @(38, "E", 43373),
    @(41, "D", 30400),
    @(41, "E", 25300),
    @(42, "D", 29200),
    @(42, "E", 41000),
    @(43, "D", 0),
    @(43, "E", 0),
    @(44, "D", 0),
    @(44, "E", 0),
    @(45, "D", 0),
    @(45, "E", 0),
    @(46, "D", 0),
    @(46, "E", 0),
    @(47, "D", 0),
    @(47, "E", 0),
    @(48, "D", 15800),
    @(48, "E", 16100),
    @(49, "D", 12400),
    @(49, "E", 12500),
    @(50, "D", 0),
    @(50, "E", 0),
    @(51, "D", 0),
    @(51, "E", 0),
    @(52, "D", 3800),
    @(52, "E", 4800),
    @(53, "D", 0),
    @(53, "E", 0),
    @(54, "D", 1455700),
    @(54, "E", 1448300),
    @(57, "D", 0),
    @(57, "E", 0),
    @(58, "D", 0),
    @(58, "E", 0),
    @(59, "D", 6500),
    @(59, "E", 6600),
    @(60, "D", 0),
    @(60, "E", 0),
    @(61, "D", 0),
    @(61, "E", 0),
    @(62, "D", 0),
    @(62, "E", 0),
    @(62, "F", 0),
    @(62, "G", 0),
    @(62, "H", 0),
    @(62, "I", 0),
    @(62, "J", 0),
    @(63, "D", 0),
    @(63, "E", 0),
    @(64, "D", 0),
    @(64, "E", 0),
    @(65, "D", 0),
    @(65, "E", 0),
    @(66, "D", 1282800),
    @(66, "E", 1279600),
    @(68, "D", 0),
    @(68, "E", 0),
    @(69, "D", 0),
    @(69, "E", 0),
    @(70, "D", 0),
    @(70, "E", 0),
    @(71, "D", 0),
    @(71, "E", 0),
    @(72, "D", 137900),
    @(72, "E", 135800),
    @(73, "D", 0),
    @(73, "E", 0),
    @(74, "D", 0),
    @(74, "E", 0),
    @(75, "D", 0),
    @(75, "E", 0),
    @(76, "D", 172900),
    @(76, "E", 168600),
    @(77, "D", 0),
    @(77, "E", 0),
    @(80, "D", 43465),
    @(80, "E", 43373),
    @(81, "D", 4200),
    @(81, "E", 4500),
    @(83, "D", 500),
    @(83, "E", 400),
    @(84, "D", 0),
    @(84, "E", 0),
    @(85, "D", 0),
    @(85, "E", 0),
    @(86, "D", 0),
    @(86, "E", 0),
    @(87, "D", 0),
    @(87, "E", 0),
    @(88, "D", 0),
    @(88, "E", 0),
    @(89, "D", 5200),
    @(89, "E", 6300),
    @(91, "D", 0),
    @(91, "E", -200),
    @(91, "H", -100),
    @(92, "D", 0),
    @(92, "E", 0),
    @(93, "D", 0),
    @(93, "E", 0),
    @(94, "D", -900),
    @(94, "E", 7700),
    @(96, "D", -2100),
    @(96, "E", -2100),
    @(97, "D", 0),
    @(97, "E", 0),
    @(98, "D", 0),
    @(98, "E", 0),
    @(99, "D", 0),
    @(99, "E", 0),
    @(100, "D", 800),
    @(100, "E", -9200),
    @(101, "D", 0),
    @(101, "E", 0),
    @(102, "D", 5100),
    @(102, "E", 4800),
)

foreach ($item in $data) {
    $r = $item[0]
    $col = $item[1]
    $v = $item[2]
    $ws.Range("$col$r").Value = $v
}
